$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 4
$ws.Range("A5").HorizontalAlignment = -4108

$ws.Range("B5").Value = "Bappaditya"
$ws.Range("E5").Value = "IOT & Robotics"
$ws.Range("F5").Value = "Inactive"
$ws.Range("C5").Value = "93398 18158"

$ws.Range("C6").Select()
